$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1901639344262295
$ws.Range("C2").Value = 0.5721311475409836
$ws.Range("J2").Value = 0.01311475409836066
$ws.Range("P2").Value = 0.140983606557377
$ws.Range("S2").Value = 0.08360655737704918
$ws.Range("B3").Value = 0.01902173913043478
$ws.Range("C3").Value = 0.04619565217391304
$ws.Range("J3").Value = 0.03260869565217391
$ws.Range("P3").Value = 0.7255434782608695
$ws.Range("S3").Value = 0.1766304347826087
$ws.Range("J4").Value = 0.02877697841726619
$ws.Range("P4").Value = 0.762589928057554
$ws.Range("S4").Value = 0.2086330935251799
$ws.Range("B6").Value = 0.07469879518072289
$ws.Range("D6").Value = 0.01445783132530121
$ws.Range("F6").Value = 0.08433734939759036
$ws.Range("J6").Value = 0.2987951807228916
$ws.Range("O6").Value = 0.009638554216867471
$ws.Range("Q6").Value = 0.1879518072289157
$ws.Range("R6").Value = 0.06024096385542169
$ws.Range("S6").Value = 0.2698795180722892
$ws.Range("B7").Value = 0.09866666666666667
$ws.Range("D7").Value = 0.04
$ws.Range("F7").Value = 0.056
$ws.Range("J7").Value = 0.136
$ws.Range("O7").Value = 0.01333333333333333
$ws.Range("Q7").Value = 0.192
$ws.Range("R7").Value = 0.064
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.1177944862155388
$ws.Range("D8").Value = 0.02255639097744361
$ws.Range("E8").Value = 0.0012531328320802
$ws.Range("F8").Value = 0.05764411027568922
$ws.Range("J8").Value = 0.1240601503759398
$ws.Range("O8").Value = 0.0137844611528822
$ws.Range("Q8").Value = 0.2205513784461153
$ws.Range("R8").Value = 0.06516290726817042
$ws.Range("S8").Value = 0.3771929824561404
$ws.Range("B9").Value = 0.1075
$ws.Range("D9").Value = 0.035
$ws.Range("F9").Value = 0.045
$ws.Range("J9").Value = 0.1125
$ws.Range("O9").Value = 0.015
$ws.Range("Q9").Value = 0.22
$ws.Range("R9").Value = 0.0625
$ws.Range("S9").Value = 0.4025
$ws.Range("B10").Value = 0.1062452399086063
$ws.Range("D10").Value = 0.03389185072353389
$ws.Range("E10").Value = 0.0007616146230007616
$ws.Range("F10").Value = 0.0594059405940594
$ws.Range("J10").Value = 0.115003808073115
$ws.Range("O10").Value = 0.01561309977151561
$ws.Range("Q10").Value = 0.2536176694592536
$ws.Range("R10").Value = 0.07730388423457731
$ws.Range("S10").Value = 0.3381568926123382
$ws.Range("G11").Value = 0.1573604060913706
$ws.Range("J11").Value = 0.07952622673434856
$ws.Range("K11").Value = 0.1979695431472081
$ws.Range("L11").Value = 0.55668358714044
$ws.Range("S11").Value = 0.008460236886632826
$ws.Range("G12").Value = 0.7138643067846607
$ws.Range("J12").Value = 0.2182890855457227
$ws.Range("K12").Value = 0.005899705014749262
$ws.Range("L12").Value = 0.02064896755162242
$ws.Range("S12").Value = 0.04129793510324484
$ws.Range("G13").Value = 0.6575342465753424
$ws.Range("J13").Value = 0.273972602739726
$ws.Range("S13").Value = 0.0684931506849315
$ws.Range("F15").Value = 0.0145985401459854
$ws.Range("H15").Value = 0.1265206812652068
$ws.Range("I15").Value = 0.07785888077858881
$ws.Range("J15").Value = 0.3892944038929441
$ws.Range("K15").Value = 0.06569343065693431
$ws.Range("M15").Value = 0.009732360097323601
$ws.Range("O15").Value = 0.06082725060827251
$ws.Range("S15").Value = 0.2554744525547445
$ws.Range("F16").Value = 0.02678571428571428
$ws.Range("H16").Value = 0.1495535714285714
$ws.Range("I16").Value = 0.08705357142857142
$ws.Range("J16").Value = 0.4330357142857143
$ws.Range("K16").Value = 0.078125
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("N16").Value = 0.004464285714285714
$ws.Range("O16").Value = 0.06473214285714286
$ws.Range("S16").Value = 0.1383928571428572
$ws.Range("F17").Value = 0.02587800369685767
$ws.Range("H17").Value = 0.1487985212569316
$ws.Range("I17").Value = 0.0933456561922366
$ws.Range("J17").Value = 0.4537892791127542
$ws.Range("K17").Value = 0.1007393715341959
$ws.Range("M17").Value = 0.01940850277264325
$ws.Range("O17").Value = 0.04990757855822551
$ws.Range("S17").Value = 0.1081330868761553
$ws.Range("F18").Value = 0.01510574018126888
$ws.Range("H18").Value = 0.1722054380664653
$ws.Range("I18").Value = 0.09969788519637462
$ws.Range("J18").Value = 0.4471299093655589
$ws.Range("K18").Value = 0.1087613293051359
$ws.Range("M18").Value = 0.01208459214501511
$ws.Range("O18").Value = 0.0513595166163142
$ws.Range("S18").Value = 0.09365558912386707
$ws.Range("F19").Value = 0.01484938481120068
$ws.Range("H19").Value = 0.1981332201951634
$ws.Range("I19").Value = 0.08230801866779805
$ws.Range("J19").Value = 0.3708103521425541
$ws.Range("K19").Value = 0.1124310564276623
$ws.Range("M19").Value = 0.01654645736105219
$ws.Range("N19").Value = 0.0008485362749257531
$ws.Range("O19").Value = 0.07297411964361476
$ws.Range("S19").Value = 0.1310988544760288
